# Swap the contents of columns E and F (codeforiati:category-name <-> codeforiati:group-code)
# for every row in the used range, including the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    # Prefix with a leading apostrophe so Excel stores the swapped value as
    # literal text instead of auto-coercing numeric-looking strings (e.g.
    # "110") into numbers. The apostrophe itself is not part of the stored
    # value.
    $eCell.Value = "'" + $fVal
    $fCell.Value = "'" + $eVal
}
